# Applies the "Add files via upload" revision to the workbook.
# Target state derived from the canonical-OOXML diff:
#  - data_path: replace rows with new source-path entries
#  - data_pattern: add ROI1/ROI2/MYROI/MEASURE columns + two cerebellum rows
#  - models: add a third model (Model_03)
#  - cosmetic: selections / column widths carried over where COM exposes them

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "demographic_clinical" — no data changes, only cosmetic tweaks
# ---------------------------------------------------------------------
$wsDemo = $wb.Worksheets.Item("demographic_clinical")
$wsDemo.Columns.Item(1).ColumnWidth = 104.6640625
$wsDemo.Range("C2").Select()

# ---------------------------------------------------------------------
# Sheet "data_path" — new list of source directories
# ---------------------------------------------------------------------
$wsPath = $wb.Worksheets.Item("data_path")

$wsPath.Range("A2").Value = "/tmp/isilon/morey/lab/dusom_morey/new_halfpipe/Outputs/preproc4d/ica"
$wsPath.Range("B2").Value = 0

$wsPath.Range("A3").Value = "/tmp/isilon/morey/lab/dusom_morey/new_halfpipe/Outputs/falff_reho"
$wsPath.Range("B3").Value = 1

$wsPath.Range("A4").Value = "/tmp/isilon/morey/lab/dusom_morey/courtney_cerebellum/enigma_Groningen/acapulco/output"
$wsPath.Range("B4").Value = 0

$wsPath.Range("A5").Value = "/tmp/isilon/morey/lab/dusom_morey/courtney_cerebellum/enigma_UMN_Lissek/acapulco/output"
$wsPath.Range("B5").Value = 0

$wsPath.Range("A6").Value = "/tmp/isilon/morey/lab/dusom_morey/courtney_cerebellum/enigma_UWMadison_Grupe/acapulco/output"
$wsPath.Range("B6").Value = 0

$wsPath.Range("A6").Select()

# ---------------------------------------------------------------------
# Sheet "data_pattern" — new ROI1 / ROI2 / MYROI / MEASURE columns
# ---------------------------------------------------------------------
$wsPattern = $wb.Worksheets.Item("data_pattern")

# Header row: E/F change meaning, G/H are new, I replaces old G (EXCLUDED)
$wsPattern.Range("E1").Value = "ROI1"
$wsPattern.Range("F1").Value = "ROI2"
$wsPattern.Range("G1").Value = "MYROI"
$wsPattern.Range("H1").Value = "MEASURE"
$wsPattern.Range("I1").Value = "EXCLUDED"

# Row 2 (fALFF_alff) — EXCLUDED flag moves from G2 to I2
$wsPattern.Range("G2").ClearContents()
$wsPattern.Range("I2").Value = 1

# Row 3 — now corrMatrix pattern (was previously different strings); EXCLUDED moves to I3
$wsPattern.Range("A3").Value = "corrMatrix_atlas-schaefer2011Combined"
$wsPattern.Range("B3").Value = "_feature-corrMatrix_atlas-schaefer2011Combined_desc-correlation_matrix.tsv"
$wsPattern.Range("C3").Value = "_feature-corrMatrix_atlas-schaefer2011Combined_timeseries.json"
$wsPattern.Range("G3").ClearContents()
$wsPattern.Range("I3").Value = 1

# Row 4 — new atlas_cerebellum_SFC entry
$wsPattern.Range("A4").Value = "atlas_cerebellum_SFC"
$wsPattern.Range("B4").Value = "_icareg_preproc_bold.nii.gz"
$wsPattern.Range("E4").Value = "/tmp/isilon/morey/lab/dusom_morey/Aurelio_preproc/atlases/tpl-MNI152NLin2009cAsym_atlas-schaefer2011Combined_dseg.nii.gz"
$wsPattern.Range("F4").Value = "_n4_mni_seg_post_inverse.nii.gz"
$wsPattern.Range("H4").Value = "Functional_Connectivity"
$wsPattern.Range("I4").Value = 1

# Row 5 — new atlas_cerebellum_SFC_MyROIs entry
$wsPattern.Range("A5").Value = "atlas_cerebellum_SFC_MyROIs"
$wsPattern.Range("B5").Value = "_icareg_preproc_bold.nii.gz"
$wsPattern.Range("E5").Value = "/tmp/isilon/morey/lab/dusom_morey/Aurelio_preproc/atlases/tpl-MNI152NLin2009cAsym_atlas-schaefer2011Combined_dseg.nii.gz"
$wsPattern.Range("F5").Value = "_n4_mni_seg_post_inverse.nii.gz"
$wsPattern.Range("G5").Value = "MY_ROIs.xlsx"
$wsPattern.Range("H5").Value = "Functional_Connectivity"
$wsPattern.Range("I5").Value = 0

$wsPattern.Columns.Item(2).ColumnWidth = 42.83203125
$wsPattern.Columns.Item(3).ColumnWidth = 62.5
$wsPattern.Columns.Item(5).ColumnWidth = 103.1640625

$wsPattern.Range("I10").Select()

# ---------------------------------------------------------------------
# Sheet "predictors" — header relabelled only (shared-string shuffle),
# no visible data change; just carry the cosmetic selection forward.
# ---------------------------------------------------------------------
$wsPred = $wb.Worksheets.Item("predictors")
$wsPred.Range("A2").Select()

# ---------------------------------------------------------------------
# Sheet "models" — add Model_03
# ---------------------------------------------------------------------
$wsModels = $wb.Worksheets.Item("models")

$wsModels.Range("A4").Value = "Model_03"
$wsModels.Range("B4").Value = "lmer(Yvar ~ GROUP * SEX + AGE +  (1|SITE))"
$wsModels.Range("C4").Value = 0

$wsModels.Range("E11").Select()
$wsModels.Activate()
